# Update the "Timothy Freeman" bio textbox on slide 2:
#   - change the degree line from "Bsc, Msc" to "BA, MPhil"
#   - rename the shape (TextBox 17 -> TextBox 16)
#   - move the shape to the end of the z-order (after "Picture 4"),
#     matching its new position in the slide's shape tree
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$shape = $s.Shapes.Item(11)

$shape.Name = "TextBox 16"

$tr = $shape.TextFrame.TextRange
[void]$tr.Replace("Bsc, Msc", "BA, MPhil")

$shape.ZOrder(0)
